# Scheduled market-data refresh: update currentAveragePrice / profit
# columns (H:N) for the affected Leve rows across the ALC, ARM, BSM,
# CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# ALC
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H76").Value = 4703.5
$ws.Range("I76").Value = 4634.5713
$ws.Range("K76").Value = 4634.5713
$ws.Range("M76").Value = -4319.5713

$ws.Range("H79").Value = 4703.5
$ws.Range("I79").Value = 4634.5713
$ws.Range("K79").Value = 4634.5713
$ws.Range("M79").Value = -3542.5713

$ws.Range("H98").Value = 472.25
$ws.Range("I98").Value = 457.72223
$ws.Range("J98").Value = 603
$ws.Range("K98").Value = 457.72223
$ws.Range("L98").Value = 603
$ws.Range("M98").Value = 1040.27777
$ws.Range("N98").Value = -3599

$ws.Range("H103").Value = 1212.1818
$ws.Range("J103").Value = 312
$ws.Range("L103").Value = 936
$ws.Range("N103").Value = -2108

$ws.Range("H122").Value = 472.25
$ws.Range("I122").Value = 457.72223
$ws.Range("J122").Value = 603
$ws.Range("K122").Value = 1373.16669
$ws.Range("L122").Value = 1809
$ws.Range("M122").Value = 1076.83331
$ws.Range("N122").Value = -6709

$ws.Range("H138").Value = 1870.8206
$ws.Range("I138").Value = 2223.3333
$ws.Range("J138").Value = 1714.1482
$ws.Range("K138").Value = 6669.999899999999
$ws.Range("L138").Value = 5142.444600000001
$ws.Range("M138").Value = -1529.999899999999
$ws.Range("N138").Value = -15422.4446

# ---------------------------------------------------------------
# ARM
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 31253
$ws.Range("I32").Value = 9119.928
$ws.Range("J32").Value = 141918.36
$ws.Range("K32").Value = 9119.928
$ws.Range("L32").Value = 141918.36
$ws.Range("M32").Value = -8832.928
$ws.Range("N32").Value = -142492.36

$ws.Range("H58").Value = 14999.667
$ws.Range("J58").Value = 14999.667
$ws.Range("L58").Value = 14999.667
$ws.Range("N58").Value = -15859.667

$ws.Range("H122").Value = 1998.8
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 2331.3333
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 6993.999899999999
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -11893.9999

# ---------------------------------------------------------------
# BSM
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H107").Value = 250003950
$ws.Range("I107").Value = 250003950
$ws.Range("K107").Value = 250003950
$ws.Range("M107").Value = -250002030

# ---------------------------------------------------------------
# CRP
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H22").Value = 417.8
$ws.Range("I22").Value = 397.25
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 397.25
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -47.25
$ws.Range("N22").Value = -1200

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

# ---------------------------------------------------------------
# CUL
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H12").Value = 40.263157
$ws.Range("I12").Value = 14.666667
$ws.Range("J12").Value = 52.076923
$ws.Range("K12").Value = 44.000001
$ws.Range("L12").Value = 156.230769
$ws.Range("M12").Value = 128.999999
$ws.Range("N12").Value = -502.230769

# ---------------------------------------------------------------
# GSM
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 107593.7
$ws.Range("I70").Value = 339069.34
$ws.Range("J70").Value = 8389.857
$ws.Range("K70").Value = 339069.34
$ws.Range("L70").Value = 8389.857
$ws.Range("M70").Value = -338799.34
$ws.Range("N70").Value = -8929.857

$ws.Range("H73").Value = 107593.7
$ws.Range("I73").Value = 339069.34
$ws.Range("J73").Value = 8389.857
$ws.Range("K73").Value = 339069.34
$ws.Range("L73").Value = 8389.857
$ws.Range("M73").Value = -338133.34
$ws.Range("N73").Value = -10261.857

$ws.Range("H122").Value = 3093.6
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3093.6
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 9280.799999999999
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -14180.8

# ---------------------------------------------------------------
# LTW
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 832.1724
$ws.Range("I22").Value = 1037.1538
$ws.Range("J22").Value = 665.625
$ws.Range("K22").Value = 1037.1538
$ws.Range("L22").Value = 665.625
$ws.Range("M22").Value = -742.1538
$ws.Range("N22").Value = -1255.625

$ws.Range("H27").Value = 832.1724
$ws.Range("I27").Value = 1037.1538
$ws.Range("J27").Value = 665.625
$ws.Range("K27").Value = 1037.1538
$ws.Range("L27").Value = 665.625
$ws.Range("M27").Value = -930.1538
$ws.Range("N27").Value = -879.625

$ws.Range("H40").Value = 42691.48
$ws.Range("I40").Value = 65355.438
$ws.Range("J40").Value = 2400
$ws.Range("K40").Value = 65355.438
$ws.Range("L40").Value = 2400
$ws.Range("M40").Value = -65219.438
$ws.Range("N40").Value = -2672

$ws.Range("H46").Value = 920736.4
$ws.Range("I46").Value = 338
$ws.Range("J46").Value = 1687735
$ws.Range("K46").Value = 338
$ws.Range("L46").Value = 1687735
$ws.Range("M46").Value = -150
$ws.Range("N46").Value = -1688111

$ws.Range("H122").Value = 4001.3333
$ws.Range("I122").Value = 4001.3333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 12003.9999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9553.999899999999
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 6211.467
$ws.Range("I132").Value = 7363
$ws.Range("J132").Value = 4895.4287
$ws.Range("K132").Value = 22089
$ws.Range("L132").Value = 14686.2861
$ws.Range("M132").Value = -19559
$ws.Range("N132").Value = -19746.2861

# ---------------------------------------------------------------
# WVR
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H122").Value = 2713.3635
$ws.Range("I122").Value = 1458.8
$ws.Range("J122").Value = 3758.8333
$ws.Range("K122").Value = 4376.4
$ws.Range("L122").Value = 11276.4999
$ws.Range("M122").Value = -1926.4
$ws.Range("N122").Value = -16176.4999

$ws.Range("H132").Value = 2443.4695
$ws.Range("I132").Value = 2662.7812
$ws.Range("J132").Value = 2030.6471
$ws.Range("K132").Value = 7988.3436
$ws.Range("L132").Value = 6091.9413
$ws.Range("M132").Value = -5458.3436
$ws.Range("N132").Value = -11151.9413

$ws.Range("H136").Value = 1604.5491
$ws.Range("I136").Value = 579.3570999999999
$ws.Range("J136").Value = 2852.6086
$ws.Range("K136").Value = 1738.0713
$ws.Range("L136").Value = 8557.825800000001
$ws.Range("M136").Value = 811.9287000000002
$ws.Range("N136").Value = -13657.8258
